$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.581.80'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '2.347.27'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'315.61"
$ws.Range("E5").Value = '  -3.53%  '
$ws.Range("D6").Value = "'107.91"
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = '  -6.51%  '
$ws.Range("D10").Value = "'41.17"
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").Value = "'0.0924"
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = "'8.50"
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = "'0.991"
$ws.Range("E14").Value = '  -5.98%  '
$ws.Range("D15").Value = "'15.89"
$ws.Range("E15").Value = '  -7.84%  '
$ws.Range("D16").Value = '2.702.03'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = '2.345.61'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '42.516.82'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = '  -4.12%  '
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("D21").Value = "'76.20"
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = "'3.58"
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = "'256.61"
$ws.Range("E23").Value = '  -8.37%  '
$ws.Range("E24").Value = '  -4.96%  '
$ws.Range("D25").Value = "'9.35"
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = "'11.35"
$ws.Range("E27").Value = '  -3.79%  '
$ws.Range("D28").Value = "'22.75"
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = '  +1.33%  '
$ws.Range("D30").Value = "'172.28"
$ws.Range("E30").Value = '  -1.32%  '
$ws.Range("D31").Value = "'36.69"
$ws.Range("E31").Value = '  -3.50%  '
$ws.Range("D32").Value = "'0.0887"
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("D33").Value = "'6.05"
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("D34").Value = "'2.86"
$ws.Range("E34").Value = '  -9.85%  '
$ws.Range("E35").Value = '  +16.50%  '
$ws.Range("D36").Value = "'0.132"
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("E37").Value = '  -6.46%  '
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("D39").Value = "'3.89"
$ws.Range("E39").Value = '  -8.64%  '
$ws.Range("D40").Value = "'2.65"
$ws.Range("E40").Value = '  -6.13%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = "'1.47"
$ws.Range("E42").Value = '  -7.22%  '
$ws.Range("D43").Value = "'70.65"
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = "'11.93"
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("D46").Value = "'111.65"
$ws.Range("E46").Value = '  -9.44%  '
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = "'85.50"
$ws.Range("E47").Value = '  -9.66%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = "'9.14"
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("D49").Value = "'5.44"
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").Value = "'74.44"
$ws.Range("E50").Value = '  +1.55%  '
$ws.Range("E51").Value = '  -2.80%  '
